# Weekly update: insert the latest week's Caqui price rows (new as-of
# date 2023-05-08) into the "Mankaki" block that already existed around
# row 24, pushing every subsequent historical row down by three.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 24; everything from old row 24 onward
# (through old row 113) shifts down to rows 27-116, matching the new
# dimension A1:T116.
$ws.Rows("24:26").Insert()

$newDate = (Get-Date -Year 2023 -Month 5 -Day 8 -Hour 0 -Minute 0 -Second 0).Date

# Row 24 - new "Especial" quote for Mankaki
$ws.Cells.Item(24, 1).Value = 8
$ws.Cells.Item(24, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(24, 3).Value = "Coquimbo"
$ws.Cells.Item(24, 4).Value = $newDate
$ws.Cells.Item(24, 5).Value = 4
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100107
$ws.Cells.Item(24, 8).Value = "Otros"
$ws.Cells.Item(24, 9).Value = 100107001
$ws.Cells.Item(24, 10).Value = "Caqui"
$ws.Cells.Item(24, 11).Value = "Mankaki"
$ws.Cells.Item(24, 12).Value = "Especial"
$ws.Cells.Item(24, 13).Value = 10
$ws.Cells.Item(24, 14).Value = 400000
$ws.Cells.Item(24, 15).Value = 410000
$ws.Cells.Item(24, 16).Value = 405000
$ws.Cells.Item(24, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(24, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(24, 19).Value = 900
$ws.Cells.Item(24, 20).Value = 450

# Row 25 - new "Primera" quote for Mankaki
$ws.Cells.Item(25, 1).Value = 8
$ws.Cells.Item(25, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(25, 3).Value = "Coquimbo"
$ws.Cells.Item(25, 4).Value = $newDate
$ws.Cells.Item(25, 5).Value = 4
$ws.Cells.Item(25, 6).Value = "Fruta"
$ws.Cells.Item(25, 7).Value = 100107
$ws.Cells.Item(25, 8).Value = "Otros"
$ws.Cells.Item(25, 9).Value = 100107001
$ws.Cells.Item(25, 10).Value = "Caqui"
$ws.Cells.Item(25, 11).Value = "Mankaki"
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 20
$ws.Cells.Item(25, 14).Value = 360000
$ws.Cells.Item(25, 15).Value = 370000
$ws.Cells.Item(25, 16).Value = 365000
$ws.Cells.Item(25, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(25, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(25, 19).Value = 811
$ws.Cells.Item(25, 20).Value = 450

# Row 26 - new "Segunda" quote for Mankaki
$ws.Cells.Item(26, 1).Value = 8
$ws.Cells.Item(26, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = $newDate
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = "Fruta"
$ws.Cells.Item(26, 7).Value = 100107
$ws.Cells.Item(26, 8).Value = "Otros"
$ws.Cells.Item(26, 9).Value = 100107001
$ws.Cells.Item(26, 10).Value = "Caqui"
$ws.Cells.Item(26, 11).Value = "Mankaki"
$ws.Cells.Item(26, 12).Value = "Segunda"
$ws.Cells.Item(26, 13).Value = 20
$ws.Cells.Item(26, 14).Value = 310000
$ws.Cells.Item(26, 15).Value = 320000
$ws.Cells.Item(26, 16).Value = 315000
$ws.Cells.Item(26, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(26, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(26, 19).Value = 700
$ws.Cells.Item(26, 20).Value = 450
